# Edit script: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# Replaces the employee debt-statement table with a new data set (9 rows instead of 4)
# and moves the legal-representative signature block further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Step 1: move the signature block (rows 24-25) down to rows 29-30, keeping
# the exact same formatting/styles, then clear out the old location.
# ---------------------------------------------------------------------------
$ws.Range("B24:C24").Copy()
$ws.Range("B29:C29").PasteSpecial($xlPasteFormats)
$ws.Range("H24:J24").Copy()
$ws.Range("H29:J29").PasteSpecial($xlPasteFormats)

$ws.Range("B25:C25").Copy()
$ws.Range("B30:C30").PasteSpecial($xlPasteFormats)
$ws.Range("H25:J25").Copy()
$ws.Range("H30:J30").PasteSpecial($xlPasteFormats)

$ws.Range("B29").Value = $ws.Range("B24").Value2
$ws.Range("H29").Value = $ws.Range("H24").Value2
$ws.Range("B30").Value = $ws.Range("B25").Value2
$ws.Range("H30").Value = $ws.Range("H25").Value2

$ws.Range("B29:C29").Merge()
$ws.Range("H29:J29").Merge()
$ws.Range("B30:C30").Merge()
$ws.Range("H30:J30").Merge()

# ---------------------------------------------------------------------------
# Step 2: grow the employee data table from 4 rows (16-19) to 9 rows (16-24).
# Row 19 used to carry the "bottom border" styling that must now live on row
# 24 (the new last row of the table); rows 16-18 carry the regular styling
# that the newly-added rows 19-23 should copy as well.
# ---------------------------------------------------------------------------
$ws.Range("B19:J19").Copy()
$ws.Range("B24:J24").PasteSpecial($xlPasteFormats)

$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial($xlPasteFormats)
$ws.Range("B20:J20").PasteSpecial($xlPasteFormats)
$ws.Range("B21:J21").PasteSpecial($xlPasteFormats)
$ws.Range("B22:J22").PasteSpecial($xlPasteFormats)
$ws.Range("B23:J23").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# Step 3: write the new employee / mora data into rows 16-24.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row = 16; Tipo = "CC"; Doc = "45484818";    Nombre = "JUDITH RUIZ LOZANO";             Periodo = "2507"; Valor = 56940;  Salario = 1423500 },
    @{ Row = 17; Tipo = "CC"; Doc = "45484818";    Nombre = "JUDITH RUIZ LOZANO";             Periodo = "2506"; Valor = 56940;  Salario = 1423500 },
    @{ Row = 18; Tipo = "CC"; Doc = "45484818";    Nombre = "JUDITH RUIZ LOZANO";             Periodo = "2505"; Valor = 18980;  Salario = 1423500 },
    @{ Row = 19; Tipo = "CC"; Doc = "1042577655";  Nombre = "CAMILA ANDREA FUENTES ARROYO";   Periodo = "2507"; Valor = 56940;  Salario = 1423500 },
    @{ Row = 20; Tipo = "CC"; Doc = "1042577655";  Nombre = "CAMILA ANDREA FUENTES ARROYO";   Periodo = "2506"; Valor = 56940;  Salario = 1423500 },
    @{ Row = 21; Tipo = "CC"; Doc = "1042577655";  Nombre = "CAMILA ANDREA FUENTES ARROYO";   Periodo = "2505"; Valor = 56940;  Salario = 1423500 },
    @{ Row = 22; Tipo = "CC"; Doc = "1042577655";  Nombre = "CAMILA ANDREA FUENTES ARROYO";   Periodo = "2504"; Valor = 56940;  Salario = 1423500 },
    @{ Row = 23; Tipo = "CC"; Doc = "45540913";    Nombre = "KEILA YOHANA BURGOS ROMERO";     Periodo = "2507"; Valor = 56940;  Salario = 1423500 },
    @{ Row = 24; Tipo = "CC"; Doc = "45540913";    Nombre = "KEILA YOHANA BURGOS ROMERO";     Periodo = "2506"; Valor = 53144;  Salario = 1423500 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.Tipo
    $ws.Cells.Item($row, 3).Value = $r.Doc
    $ws.Cells.Item($row, 4).Value = $r.Nombre
    $ws.Cells.Item($row, 5).Value = $r.Periodo
    $ws.Cells.Item($row, 6).Value = $r.Valor
    $ws.Cells.Item($row, 7).Value = $r.Salario
}

# ---------------------------------------------------------------------------
# Step 4: misc value updates elsewhere on the sheet.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 470704
$ws.Range("F13").Value = 4

Write-Host "edit applied"
